$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$tbl = $ws.ListObjects.Item(1)

# Extend the "Tabela1" table by one row (mirrors Excel's table auto-expand
# when data is typed directly below it) and populate the new row.
$newRow = $tbl.ListRows.Add()
$ws.Range("A45").Value = "ArrDelayGroup"
$ws.Range("B45").Value = 43

# The previous last row (44) keeps its existing "bottom of table" border
# look; only the brand-new last row needs formatting. Column A of the new
# row takes the regular body-row look (same as A43), while column B keeps
# the special bottom-border look the previous last row (B44) had.
$ws.Range("A43").Copy()
$ws.Range("A45").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B44").Copy()
$ws.Range("B45").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Move the selection to the newly added cell, as in the authored workbook.
$ws.Range("A45").Select()
